# Update the EPEX Spot prices workbook with the latest day of data.

$wb = $excel.ActiveWorkbook

# --- Sheet "Prix Spot": add a new day column (G) for 20-jun ---
$wsSpot = $wb.Worksheets.Item("Prix Spot")

$wsSpot.Range("F1").Copy()
$wsSpot.Range("G1").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$wsSpot.Range("G1").Value = "20-jun"

$spotValues = @(55.73, 21.73, 25.89, 25.45, 13.47, 30.38, 37.45, 57.51, 67.17, 38.01, 4.48, 0, -0.01, -0.01, -0.01, 0.05, 12.8, 50.91, 94.78, 115, 126.68, 115.67, 122.86, 114.9)

for ($i = 0; $i -lt $spotValues.Length; $i++) {
    $row = $i + 2
    $wsSpot.Cells.Item($row, 7).Value = $spotValues[$i]
}

# --- Sheet "Gaz": append row for 2025-06-18 ---
$wsGaz = $wb.Worksheets.Item("Gaz")
# Force the date to be stored as plain text (matching the existing rows)
# instead of letting Excel auto-convert it to a date serial number.
$wsGaz.Cells.Item(4, 1).NumberFormat = "@"
$wsGaz.Cells.Item(4, 1).Value = "2025-06-18"
$wsGaz.Cells.Item(4, 1).ClearFormats()
$wsGaz.Cells.Item(4, 2).Value = 38.45
$wsGaz.Cells.Item(4, 3).Value = 10800
$wsGaz.Cells.Item(4, 4).Value = 37.806

# --- Sheet "CO2": append row for 2025-06-18 ---
$wsCo2 = $wb.Worksheets.Item("CO2")
$wsCo2.Cells.Item(4, 1).NumberFormat = "@"
$wsCo2.Cells.Item(4, 1).Value = "2025-06-18"
$wsCo2.Cells.Item(4, 1).ClearFormats()
$wsCo2.Cells.Item(4, 2).Value = 73.45
